# Mimic importing rows from data.csv:
#  - a date value (formatted DD/MM/YY)
#  - a blank cell formatted as a US-style date ([$-409]MM/DD/YY)
#  - drop the old "a" label
#  - two new rows of shared-string labels (z/x/c)
#  - a TODAY() formula (formatted DD/MM/YY) plus an arrow glyph label

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: turn A1 into a date value (serial 34262 = 1993-10-20), clear B1 but
# give it a date number format.
$ws.Range("A1").Value = 34262
$ws.Range("A1").NumberFormat = "DD/MM/YY"

$ws.Range("B1").Value = $null
$ws.Range("B1").NumberFormat = "[$-409]MM/DD/YY"

# Row 2: drop the old "a" label from A2, keep "b"/"c" in B2/C2.
$ws.Range("A2").ClearContents()

# Row 4: new row of shared-string labels, mimicking an extra CSV row.
$ws.Range("A4").Value = "z"
$ws.Range("B4").Value = "x"
$ws.Range("C4").Value = "c"

# Row 6: a volatile TODAY() formula plus an arrow glyph label.
# (Number format is applied before the formula so the engine doesn't stamp
# a transient default date format into the style table first.)
$ws.Range("A6").NumberFormat = "DD/MM/YY"
$ws.Range("A6").Formula = "=TODAY()"
$ws.Range("C6").Value = [char]8592

# Move the active selection like the end of the editing session.
$ws.Range("B10").Select() | Out-Null
